# Fruta / hortaliza, semanal
# The commit reshuffles the data rows (rows 2-40) of the single sheet: every
# row's full set of column values (A:R) moves to a different row position.
# No values are actually changed/added/removed - the 39 data rows are simply
# permuted. Column D holds a date (stored as the Excel serial number), so we
# read/write it through .Value2 to avoid locale/COM date coercion issues.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 40
$lastCol = 18   # A..R

# Map: new row number -> row number that currently holds the data that must
# end up there (1-based data rows, matching the sheet's own row numbers).
$rowMap = @{
    2  = 11
    3  = 3
    4  = 23
    5  = 31
    6  = 4
    7  = 37
    8  = 14
    9  = 25
    10 = 16
    11 = 2
    12 = 30
    13 = 35
    14 = 34
    15 = 26
    16 = 12
    17 = 9
    18 = 40
    19 = 21
    20 = 22
    21 = 10
    22 = 28
    23 = 7
    24 = 5
    25 = 39
    26 = 24
    27 = 27
    28 = 33
    29 = 6
    30 = 8
    31 = 15
    32 = 17
    33 = 18
    34 = 36
    35 = 32
    36 = 19
    37 = 20
    38 = 13
    39 = 38
    40 = 29
}

# Snapshot every existing row's values (A:R) before writing anything, since
# several rows are sources for more than one destination and vice versa.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the permuted data back out.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowMap[$r]
    $rowVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$c]
    }
}
